# Slide 5 ("Problemas Reais que o Projeto visa Resolver"): drop the
# " - Ricardo" author suffix from the title and fill in the body
# placeholder with the list of problems the project tries to solve.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

# --- Title -----------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Problemas Reais que o Projeto visa Resolver"

# --- Body placeholder --------------------------------------------------
$body = $s.Shapes.Item(2)

# Resize/position the placeholder (matches the manual resize done while
# authoring the bullet list) and let the text shrink to fit.
$body.Left = 90.94125747680665
$body.Top = 204.99991607666018
$body.Width = 694.9338073730469
$body.Height = 311.18464660644537
$body.TextFrame.AutoSize = 2

$intro = "O Projectary visa resolver os seguintes problemas:"
$item1 = "Inexistência de uma plataforma onde os alunos possa ter conhecimento de todos os projetos finais aos quais se podem candidatar;"
$item2 = "Ausência de uma plataforma onde os alunos se podem candidatar aos projetos finais;"
$item3 = "Ausência de um sítio onde os professores podem dispor informações aos alunos e onde os alunos podem ter acesso a essas informações;"
$item4 = "Inexistência de uma plataforma onde todas as pessoas podem ter acesso aos projetos finais já realizados e a todas as informações que foram necessárias para a realização desse projeto;"

# Seed the placeholder with a single run/paragraph and set the language
# on it so later re-typed paragraphs inherit pt-PT too.
$body.TextFrame.TextRange.Text = $intro
$body.TextFrame.TextRange.LanguageID = "pt-PT"

$body.TextFrame.TextRange.Text = $intro + "`r" + $item1 + "`r" + $item2 + "`r" + $item3 + "`r" + $item4

$tr = $body.TextFrame.TextRange

# Second-level bullets for the four problem statements.
$tr.Paragraphs(2).IndentLevel = 2
$tr.Paragraphs(3).IndentLevel = 2
$tr.Paragraphs(4).IndentLevel = 2
$tr.Paragraphs(5).IndentLevel = 2

# Italicise the product name within the intro sentence.
$nameStart = $intro.IndexOf("Projectary") + 1
$tr.Characters($nameStart, 10).Font.Italic = $true
